$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.37
$ws.Range("D3").Value = 0.37
$ws.Range("E2").Value = 0.124
$ws.Range("E3").Value = 0.124
$ws.Range("F2").Value = 0.146
$ws.Range("F3").Value = 0.146
$ws.Range("G2").Value = 0.3767483870967742
$ws.Range("G3").Value = 0.3767483870967742
$ws.Range("H2").Value = 0.3767483870967742
$ws.Range("H3").Value = 0.3767483870967742
$ws.Range("I2").Value = 0.3855612903225806
$ws.Range("I3").Value = 0.3855612903225806
$ws.Range("J2").Value = 0.2868127681012968
$ws.Range("J3").Value = 0.2868127681012968
$ws.Range("K2").Value = 478.9
$ws.Range("K3").Value = 478.9
$ws.Range("L2").Value = 0.06179354838709677
$ws.Range("L3").Value = 0.06179354838709677
$ws.Range("M2").Value = 42.2
$ws.Range("M3").Value = 42.2
$ws.Range("N2").Value = 0.002174249059714565
$ws.Range("N3").Value = 0.002174249059714565
$ws.Range("O2").Value = 0.08811860513677178
$ws.Range("O3").Value = 0.08811860513677178
$ws.Range("P2").Value = 42.2
$ws.Range("P3").Value = 42.2
$ws.Range("Q2").Value = 0.002174249059714565
$ws.Range("Q3").Value = 0.002174249059714565
$ws.Range("R2").Value = 0.08811860513677178
$ws.Range("R3").Value = 0.08811860513677178
$ws.Range("U2").Value = 417.2
$ws.Range("U3").Value = 417.2
$ws.Range("V2").Value = 0.02149518264722551
$ws.Range("V3").Value = 0.02149518264722551
$ws.Range("W2").Value = 0.1298641429617377
$ws.Range("W3").Value = 0.1298641429617377
$ws.Range("X2").Value = 0.07701817426949743
$ws.Range("X3").Value = 0.07701817426949743
$ws.Range("Y2").Value = 0.05284596869224024
$ws.Range("Y3").Value = 0.05284596869224024
$ws.Range("Z2").Value = 0.4416583558911526
$ws.Range("Z3").Value = 0.4416583558911526
$ws.Range("AA2").Value = 0.1266732556082091
$ws.Range("AA3").Value = 0.1266732556082091
$ws.Range("AB2").Value = 0.05822969170412053
$ws.Range("AB3").Value = 0.05822969170412053
$ws.Range("AC2").Value = 0.06844356390408861
$ws.Range("AC3").Value = 0.06844356390408861
$ws.Range("AD2").Value = 13733.7
$ws.Range("AD3").Value = 13733.7
$ws.Range("AF2").Value = 13733.7
$ws.Range("AF3").Value = 13733.7
$ws.Range("AG2").Value = 13316.5
$ws.Range("AG3").Value = 13316.5
$ws.Range("AH2").Value = 0.4143808440471055
$ws.Range("AH3").Value = 0.4143808440471055
$ws.Range("AI2").Value = 0.6486452021215799
$ws.Range("AI3").Value = 0.6486452021215799
$ws.Range("AJ2").Value = 0.406915096789965
$ws.Range("AJ3").Value = 0.406915096789965
$ws.Range("AK2").Value = 0.6415827941240237
$ws.Range("AK3").Value = 0.6415827941240237
$ws.Range("AL2").Value = 1320.2
$ws.Range("AL3").Value = 1320.2
$ws.Range("AM2").Value = 1320.2
$ws.Range("AM3").Value = 1320.2
$ws.Range("AN2").Value = 4.543971678136581
$ws.Range("AN3").Value = 4.543971678136581
$ws.Range("AO2").Value = 2.263369186486896
$ws.Range("AO3").Value = 2.263369186486896
$ws.Range("AP2").Value = 4.405935680254102
$ws.Range("AP3").Value = 4.405935680254102
$ws.Range("AQ2").Value = 2.263369186486896
$ws.Range("AQ3").Value = 2.263369186486896
